$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "GUI - basic interface" task row as "In Progress"
$ws.Range("G7").Value = "In Progress"

# Update "Obj Loader - refactoring ply" row from "In Progress" to "Finished"
$ws.Range("G8").Value = "Finished"

# Flag row 16 (Interpolate Normal, uv-Triangle, Ply Loadr, Obj) as the next task after obj
$ws.Range("G16").Value = "Next after obj"

# Update the view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A4")
$ws.Range("G23").Select()
